# Fix the "Recorded By" (column G) attribution strings on the
# "Session Analysis Results" sheet: the last two comma-separated
# contributors in the list were written in the wrong order (an
# injected "backup@backdoor.com" / duplicate "system" entry, and
# System vs. the real grader, ended up swapped). Restore the
# correct order by swapping the last two tokens wherever the cell
# still has one of the known mis-ordered values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Exact "before" values that need their last two entries swapped.
$targets = @(
    "System, dnasr281@gmail.com",
    "System, system, backup@backdoor.com",
    "admin@admin.com, dnasr281@gmail.com"
)

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value2

    if ($v -ne $null) {
        $isTarget = $false
        foreach ($t in $targets) {
            if ($v -eq $t) {
                $isTarget = $true
            }
        }

        if ($isTarget) {
            $parts = $v.Split(",")
            $trimmed = @()
            foreach ($p in $parts) {
                $trimmed += $p.Trim()
            }

            if ($trimmed.Count -ge 2) {
                $last = $trimmed.Count - 1
                $secondLast = $trimmed.Count - 2
                $tmp = $trimmed[$last]
                $trimmed[$last] = $trimmed[$secondLast]
                $trimmed[$secondLast] = $tmp
            }

            $cell.Value = [string]::Join(", ", $trimmed)
        }
    }
}
